$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.762.30'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').Value = '  +0.31%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.674.15'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').Value = '  +0.99%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.86'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.29'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').Value = '  +0.61%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').Value = '  +5.67%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.123'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').Value = '  +0.25%  '

$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('E11').Value = '  -0.13%  '

$ws.Range('E12').Value = '  -0.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.50'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').Value = '  -0.74%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000197'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.154.82'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.597.53'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').Value = '  +0.46%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.671.67'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').Value = '  +0.88%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.68'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').Value = '  +0.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.83'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').Value = '  -1.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.57'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').Value = '  +1.63%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '352.17'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').Value = '  -1.83%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.94'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').Value = '  +0.27%  '

$ws.Range('E24').Value = '  +5.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.81'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').Value = '  +4.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.63'
$ws.Range('D26').ClearFormats()

$ws.Range('E26').Value = '  -4.28%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.169'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').Value = '  +2.21%  '

$ws.Range('E28').Value = '  -1.47%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.13'
$ws.Range('D29').ClearFormats()

$ws.Range('E29').Value = '  +0.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '543.19'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').Value = '  +3.61%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.16'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').Value = '  -0.79%  '

$ws.Range('E33').Value = '  +0.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.59'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').Value = '  +4.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.48'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').Value = '  -0.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.424'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').Value = '  -1.61%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.44'
$ws.Range('D37').ClearFormats()

$ws.Range('E37').Value = '  -0.91%  '

$ws.Range('E38').Value = '  +0.03%  '

$ws.Range('B39').Value = 'Monero'

$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '158.38'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').Value = '  -2.18%  '

$ws.Range('B40').Value = 'Stacks'

$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.95'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').Value = '  -1.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.90'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').Value = '  +2.61%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '165.39'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').Value = '  -0.08%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.08'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').Value = '  -1.06%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0614'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').Value = '  +1.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').Value = '  -2.08%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.31'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').Value = '  +1.66%  '

$ws.Range('E48').Value = '  -0.86%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0259'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').Value = '  -1.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.102'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').Value = '  +3.82%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.30'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').Value = '  +3.38%  '
